# Fruta / hortaliza, semanal
# Insert a new weekly record as row 30 (shifting the existing rows 30-88 down
# to 31-89), matching the OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data down one row, starting at row 30.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly record.
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 44935
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103003
$ws.Range("J30").Value = "Damasco"
$ws.Range("K30").Value = "Modesto"
$ws.Range("L30").Value = "Especial"
$ws.Range("M30").Value = 500
$ws.Range("N30").Value = 18000
$ws.Range("O30").Value = 18000
$ws.Range("P30").Value = 18000
$ws.Range("Q30").Value = "$/bandeja 10 kilos"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 1800
$ws.Range("T30").Value = 10
